$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.059.94'
$ws.Range('E2').Value = '  -0.35%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.421.56'
$ws.Range('E3').Value = '  -0.18%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.44'
$ws.Range('E5').Value = '  -0.15%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.76'
$ws.Range('E6').Value = '  -0.63%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.530'
$ws.Range('E8').Value = '  -0.53%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.420.68'
$ws.Range('E9').Value = '  -0.16%  '

# Row 10
$ws.Range('E10').Value = '  -0.73%  '

# Row 11
$ws.Range('E11').Value = '  +0.22%  '

# Row 12
$ws.Range('E12').Value = '  -3.39%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.350'
$ws.Range('E13').Value = '  -1.12%  '

# Row 14
$ws.Range('E14').Value = '  +0.75%  '

# Row 15
$ws.Range('E15').Value = '  -2.54%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.859.62'
$ws.Range('E16').Value = '  -0.12%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.923.52'
$ws.Range('E17').Value = '  -0.33%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.427.60'
$ws.Range('E18').Value = '  +0.08%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.25'
$ws.Range('E19').Value = '  -1.05%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '323.17'
$ws.Range('E20').Value = '  -0.46%  '

# Row 21
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.14'
$ws.Range('E21').Value = '  -1.64%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.82'
$ws.Range('E22').Value = '  +0.66%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.10%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.27'
$ws.Range('E24').Value = '  +2.54%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.73'
$ws.Range('E25').Value = '  +0.31%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.77'
$ws.Range('E26').Value = '  -1.75%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '557.37'
$ws.Range('E27').Value = '  -5.46%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.541.27'
$ws.Range('E28').Value = '  +0.50%  '

# Row 29
$ws.Range('E29').Value = '  -0.07%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0932'
$ws.Range('E30').Value = '  -1.89%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.23'
$ws.Range('E31').Value = '  -0.75%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.39'
$ws.Range('E32').Value = '  -5.93%  '

# Row 33
$ws.Range('E33').Value = '  -1.87%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.87'
$ws.Range('E34').Value = '  -0.97%  '

# Row 35
$ws.Range('E35').Value = '  -3.58%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.02%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.74'
$ws.Range('E37').Value = '  -1.73%  '

# Row 38
$ws.Range('E38').Value = '  -1.16%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.46'
$ws.Range('E39').Value = '  -5.33%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '152.12'
$ws.Range('E40').Value = '  -1.01%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.66'
$ws.Range('E41').Value = '  -0.19%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.80'
$ws.Range('E42').Value = '  -2.09%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.17%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.26'
$ws.Range('E44').Value = '  -3.79%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '147.35'
$ws.Range('E45').Value = '  -2.15%  '

# Row 46
$ws.Range('E46').Value = '  -0.94%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0532'
$ws.Range('E47').Value = '  -1.65%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.93'
$ws.Range('E48').Value = '  -2.79%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.594'
$ws.Range('E49').Value = '  -0.08%  '

# Row 50
$ws.Range('E50').Value = '  -0.48%  '

# Row 51
$ws.Range('E51').Value = '  -0.35%  '
